$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.937.84'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '3.372.80'
$ws.Range('E3').Value = '  -3.20%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.65'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.44'
$ws.Range('E6').Value = '  +3.47%  '
$ws.Range('E7').Value = '  +5.02%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '3.370.48'
$ws.Range('E9').Value = '  -3.15%  '
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('D13').Value = '3.960.53'
$ws.Range('E13').Value = '  -2.99%  '
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.94'
$ws.Range('E15').Value = '  -3.79%  '
$ws.Range('D16').Value = '65.899.90'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').Value = '3.377.69'
$ws.Range('E18').Value = '  -2.95%  '
$ws.Range('E19').Value = '  -2.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.60'
$ws.Range('E20').Value = '  -2.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.58'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('E22').Value = '  -3.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.50'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.73'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.98'
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.73'
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '23.06'
$ws.Range('E32').Value = '  -4.32%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  -2.38%  '
$ws.Range('E35').Value = '  -4.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.51'
$ws.Range('E36').Value = '  -2.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.86'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('E38').Value = '  -3.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.05'
$ws.Range('E39').Value = '  -8.53%  '
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('D42').Value = '2.685.28'
$ws.Range('E42').Value = '  -4.50%  '
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('E44').Value = '  -3.28%  '
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '338.17'
$ws.Range('E46').Value = '  +9.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.84'
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.27'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('E50').Value = '  +2.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.28'
$ws.Range('E51').Value = '  +2.39%  '
